# Generate Report for Handoff
# Refresh "Latest Handoff Date(time)" stamps for rows that are currently
# "Handback transform failed" (row 7) or "Ready for handoff" (rows 10-16)
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$ovRows = @(7, 10, 11, 12, 13, 14, 15, 16)

$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $ovRows) {
    $wsOverview.Cells.Item($r, 4).Value = "2016-03-21 16:32:11"
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $ovRows) {
    $wsZhCn.Cells.Item($r, 5).Value = "2016-03-21 16:32:07"
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $ovRows) {
    $wsDeDe.Cells.Item($r, 5).Value = "2016-03-21 16:32:11"
}
